# Commit: "Added R and S measurement pictures"
#
# The chart ("Chart 2") that plots BBY40 Capacitance Vs. Reverse Voltage was
# moved/resized on Sheet1 to make room for newly-added R and S measurement
# pictures, and the active cell selection moved from V15 to H24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move/resize the existing chart from its old spot (roughly H3:P17) down to
# A17:H31 to make room for the new pictures above/beside it.
# Values are expressed in points, matching the target <xdr:twoCellAnchor>
# from col=0/colOff=423862, row=16/rowOff=0 to col=7/colOff=423862,
# row=30/rowOff=76200 (EMU_PER_PT = 12700; row height = 15pt; col widths
# taken from the sheet's actual column widths).
$co = $ws.ChartObjects().Item(1)
$co.Left = 33.37496062992126
$co.Top = 240.0
$co.Width = 426.357421875
$co.Height = 216.0

# Update the saved selection/active cell on Sheet1.
$ws.Range("H24").Select()
